$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cell H1 (bold, centered, bordered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I and J, rows 2-31
$data = @{
    2  = @(3, 3)
    3  = @(1, 5)
    4  = @(1, 3)
    5  = @(1, 5)
    6  = @(1, 5)
    7  = @(1, 5)
    8  = @(1, 5)
    9  = @(1, 6)
    10 = @(1, 5)
    11 = @(1, 6)
    12 = @(1, 7)
    13 = @(1, 6)
    14 = @(1, 6)
    15 = @(1, 5)
    16 = @(1, 4)
    17 = @(1, 6)
    18 = @(1, 6)
    19 = @(1, 6)
    20 = @(1, 5)
    21 = @(1, 3)
    22 = @(1, 4)
    23 = @(1, 5)
    24 = @(1, 5)
    25 = @(1, 6)
    26 = @(1, 4)
    27 = @(1, 4)
    28 = @(1, 4)
    29 = @(3, 3)
    30 = @(3, 3)
    31 = @(1, 1)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}

# Row 31's H column value is unchanged (already 1), but ensure it stays set.
$ws.Range("H31").Value = 1
